$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = "alice.smith@email.com"
$ws.Range("I3").Select()
